# Updates the C3DC phs002790 "Response-NotDone" queries so the JOINs use
# the renamed id columns (id -> study_id / participant_id) that now match
# the source dataframes, across every StatQuery / TabQuery cell on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldJoin = "LEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"`nLEFT JOIN `n    df_treatments trt ON prt.id = trt.`"participant.id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.id = trr.`"participant.id`"`nLEFT JOIN `n    df_survival srv ON prt.id = srv.`"participant.id`"`nLEFT JOIN `n    df_reference_files rfs ON std.id = rfs.`"study.id`""

$newJoin = "LEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`nLEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`nLEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`nLEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`""

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $text = $cell.Value2
        if ($text -ne $null -and $text -is [string] -and $text.Contains($oldJoin)) {
            $cell.Value2 = $text.Replace($oldJoin, $newJoin)
        }
    }
}

# Column C was resized by hand (no longer an auto "best fit" width).
$ws.Columns("C").ColumnWidth = 72.1640625
